$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F, matching the style of existing header cells (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Populate time_taken values for each data row (F2:F69)
$values = @(
    "2021-10-05 13:40:32.017027",
    "2021-10-05 13:40:32.017040",
    "2021-10-05 13:40:32.017044",
    "2021-10-05 13:40:32.017047",
    "2021-10-05 13:40:32.017050",
    "2021-10-05 13:40:32.017053",
    "2021-10-05 13:40:32.017056",
    "2021-10-05 13:40:32.017059",
    "2021-10-05 13:40:32.017063",
    "2021-10-05 13:40:32.017066",
    "2021-10-05 13:40:32.017069",
    "2021-10-05 13:40:32.017071",
    "2021-10-05 13:40:32.017075",
    "2021-10-05 13:40:32.017077",
    "2021-10-05 13:40:32.017091",
    "2021-10-05 13:40:32.017095",
    "2021-10-05 13:40:32.017098",
    "2021-10-05 13:40:32.017101",
    "2021-10-05 13:40:32.017104",
    "2021-10-05 13:40:32.017107",
    "2021-10-05 13:40:32.017110",
    "2021-10-05 13:40:32.017113",
    "2021-10-05 13:40:32.017116",
    "2021-10-05 13:40:32.017119",
    "2021-10-05 13:40:32.017122",
    "2021-10-05 13:40:32.017125",
    "2021-10-05 13:40:32.017128",
    "2021-10-05 13:40:32.017131",
    "2021-10-05 13:40:32.017134",
    "2021-10-05 13:40:32.017137",
    "2021-10-05 13:40:32.017140",
    "2021-10-05 13:40:32.017143",
    "2021-10-05 13:40:32.017146",
    "2021-10-05 13:40:32.017150",
    "2021-10-05 13:40:32.017153",
    "2021-10-05 13:40:32.017156",
    "2021-10-05 13:40:32.017159",
    "2021-10-05 13:40:32.017162",
    "2021-10-05 13:40:32.017165",
    "2021-10-05 13:40:32.017168",
    "2021-10-05 13:40:32.017171",
    "2021-10-05 13:40:32.017174",
    "2021-10-05 13:40:32.017177",
    "2021-10-05 13:40:32.017180",
    "2021-10-05 13:40:32.017183",
    "2021-10-05 13:40:32.017186",
    "2021-10-05 13:40:32.017189",
    "2021-10-05 13:40:32.017192",
    "2021-10-05 13:40:32.017195",
    "2021-10-05 13:40:32.017198",
    "2021-10-05 13:40:32.017201",
    "2021-10-05 13:40:32.017204",
    "2021-10-05 13:40:32.017207",
    "2021-10-05 13:40:32.017211",
    "2021-10-05 13:40:32.017214",
    "2021-10-05 13:40:32.017217",
    "2021-10-05 13:40:32.017220",
    "2021-10-05 13:40:32.017223",
    "2021-10-05 13:40:32.017226",
    "2021-10-05 13:40:32.017228",
    "2021-10-05 13:40:32.017232",
    "2021-10-05 13:40:32.017235",
    "2021-10-05 13:40:32.017238",
    "2021-10-05 13:40:32.017241",
    "2021-10-05 13:40:32.017245",
    "2021-10-05 13:40:32.017249",
    "2021-10-05 13:40:32.017252",
    "2021-10-05 13:40:32.017255"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
